$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3078
$ws.Range("J40").Value = 3159.8
$ws.Range("L40").Value = 3159.8
$ws.Range("N40").Value = -3509.8

$ws.Range("H88").Value = 5999.909
$ws.Range("J88").Value = 6409.9
$ws.Range("L88").Value = 6409.9
$ws.Range("N88").Value = -7221.9

$ws.Range("H91").Value = 5999.909
$ws.Range("J91").Value = 6409.9
$ws.Range("L91").Value = 6409.9
$ws.Range("N91").Value = -9217.9

$ws.Range("H94").Value = 4717.1665
$ws.Range("I94").Value = 3695.6
$ws.Range("K94").Value = 3695.6
$ws.Range("M94").Value = -3244.6

$ws.Range("H98").Value = 1052.2941
$ws.Range("I98").Value = 1020.7857
$ws.Range("K98").Value = 1020.7857
$ws.Range("M98").Value = 477.2143

$ws.Range("H100").Value = 4394.5884
$ws.Range("I100").Value = 1839.25
$ws.Range("K100").Value = 1839.25
$ws.Range("M100").Value = -1298.25

$ws.Range("H122").Value = 1052.2941
$ws.Range("I122").Value = 1020.7857
$ws.Range("K122").Value = 3062.3571
$ws.Range("M122").Value = -612.3571000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 8000
$ws.Range("I6").Value = 15000
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = -14827
$ws.Range("N6").Value = -1346

$ws.Range("H32").Value = 2715.1187
$ws.Range("I32").Value = 2268.5818
$ws.Range("K32").Value = 2268.5818
$ws.Range("M32").Value = -1981.5818

$ws.Range("H45").Value = 6199.6
$ws.Range("I45").Value = 5999.5
$ws.Range("K45").Value = 5999.5
$ws.Range("M45").Value = -5622.5

$ws.Range("H74").Value = 1755.7727
$ws.Range("I74").Value = 1601.5714
$ws.Range("K74").Value = 1601.5714
$ws.Range("M74").Value = -727.5714

$ws.Range("H77").Value = 1755.7727
$ws.Range("I77").Value = 1601.5714
$ws.Range("K77").Value = 8007.857
$ws.Range("M77").Value = -3639.857

$ws.Range("H97").Value = 680.2353000000001
$ws.Range("I97").Value = 598.4545000000001
$ws.Range("J97").Value = 830.1667
$ws.Range("K97").Value = 598.4545000000001
$ws.Range("L97").Value = 830.1667
$ws.Range("M97").Value = -102.4545000000001
$ws.Range("N97").Value = -1822.1667

$ws.Range("H122").Value = 3125.2778
$ws.Range("I122").Value = 2819.7693
$ws.Range("K122").Value = 8459.3079
$ws.Range("M122").Value = -6009.3079

$ws.Range("H132").Value = 1834.1364
$ws.Range("I132").Value = 1868.8948
$ws.Range("K132").Value = 5606.6844
$ws.Range("M132").Value = -3076.6844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 5000
$ws.Range("I14").Value = 5000
$ws.Range("K14").Value = 5000
$ws.Range("M14").Value = -4828

$ws.Range("H86").Value = 3005.5557
$ws.Range("I86").Value = 2230.8462
$ws.Range("K86").Value = 2230.8462
$ws.Range("M86").Value = -1107.8462

$ws.Range("H89").Value = 3005.5557
$ws.Range("I89").Value = 2230.8462
$ws.Range("K89").Value = 11154.231
$ws.Range("M89").Value = -5538.231

$ws.Range("H94").Value = 794.25
$ws.Range("I94").Value = 799.3200000000001
$ws.Range("K94").Value = 799.3200000000001
$ws.Range("M94").Value = -348.3200000000001

$ws.Range("H105").Value = 5045.615
$ws.Range("I105").Value = 4308.3
$ws.Range("J105").Value = 7503.3335
$ws.Range("K105").Value = 4308.3
$ws.Range("L105").Value = 7503.3335
$ws.Range("M105").Value = -2561.3
$ws.Range("N105").Value = -10997.3335

$ws.Range("H107").Value = 6333.3125
$ws.Range("I107").Value = 4666.636
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 4666.636
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -2746.636
$ws.Range("N107").Value = -13840

$ws.Range("H134").Value = 2178.1724
$ws.Range("I134").Value = 2083.577
$ws.Range("J134").Value = 2998
$ws.Range("K134").Value = 6250.731000000001
$ws.Range("L134").Value = 8994
$ws.Range("M134").Value = -3715.731000000001
$ws.Range("N134").Value = -14064

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3058.35
$ws.Range("J31").Value = 3984
$ws.Range("L31").Value = 3984
$ws.Range("N31").Value = -4574

$ws.Range("H34").Value = 3058.35
$ws.Range("J34").Value = 3984
$ws.Range("L34").Value = 3984
$ws.Range("N34").Value = -4388

$ws.Range("H105").Value = 2470.4614
$ws.Range("I105").Value = 2819.5
$ws.Range("J105").Value = 2171.2856
$ws.Range("K105").Value = 2819.5
$ws.Range("L105").Value = 2171.2856
$ws.Range("M105").Value = -1072.5
$ws.Range("N105").Value = -5665.2856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1199
$ws.Range("I51").Value = 797.6667
$ws.Range("J51").Value = 1500
$ws.Range("K51").Value = 2393.0001
$ws.Range("L51").Value = 4500
$ws.Range("M51").Value = -1933.0001
$ws.Range("N51").Value = -5420

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H98").Value = 149
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 149
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 447
$ws.Range("N98").Value = -3443
$ws.Range("M98").ClearContents()

$ws.Range("H107").Value = 1523.625
$ws.Range("I107").Value = 2680.25
$ws.Range("K107").Value = 8040.75
$ws.Range("M107").Value = -6120.75

$ws.Range("H113").Value = 2897.4285
$ws.Range("I113").Value = 945
$ws.Range("K113").Value = 2835
$ws.Range("M113").Value = -665

$ws.Range("H124").Value = 685
$ws.Range("I124").Value = 685
$ws.Range("K124").Value = 2055
$ws.Range("M124").Value = 2855

$ws.Range("H131").Value = 1844.4595
$ws.Range("J131").Value = 2113.3572
$ws.Range("L131").Value = 6340.071599999999
$ws.Range("N131").Value = -16420.0716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H107").Value = 553.04346
$ws.Range("I107").Value = 655.5833
$ws.Range("J107").Value = 441.18182
$ws.Range("K107").Value = 655.5833
$ws.Range("L107").Value = 441.18182
$ws.Range("M107").Value = 1264.4167
$ws.Range("N107").Value = -4281.18182

$ws.Range("H132").Value = 1982.08
$ws.Range("I132").Value = 1787.8096
$ws.Range("K132").Value = 5363.4288
$ws.Range("M132").Value = -2833.4288

$ws.Range("H136").Value = 28999.8
$ws.Range("J136").Value = 28999.8
$ws.Range("L136").Value = 86999.39999999999
$ws.Range("N136").Value = -92099.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4344.074
$ws.Range("I61").Value = 3046.842
$ws.Range("J61").Value = 7425
$ws.Range("K61").Value = 3046.842
$ws.Range("L61").Value = 7425
$ws.Range("M61").Value = -2844.842
$ws.Range("N61").Value = -7829

$ws.Range("H93").Value = 3124.9443
$ws.Range("I93").Value = 614.3182
$ws.Range("K93").Value = 614.3182
$ws.Range("M93").Value = 633.6818

$ws.Range("H113").Value = 4344.074
$ws.Range("I113").Value = 3046.842
$ws.Range("J113").Value = 7425
$ws.Range("K113").Value = 3046.842
$ws.Range("L113").Value = 7425
$ws.Range("M113").Value = -876.8420000000001
$ws.Range("N113").Value = -11765

$ws.Range("H122").Value = 7211.757
$ws.Range("I122").Value = 6093.04
$ws.Range("K122").Value = 18279.12
$ws.Range("M122").Value = -15829.12

$ws.Range("H136").Value = 23988.174
$ws.Range("I136").Value = 1431.1666
$ws.Range("J136").Value = 48595.816
$ws.Range("K136").Value = 4293.4998
$ws.Range("L136").Value = 145787.448
$ws.Range("M136").Value = -1743.4998
$ws.Range("N136").Value = -150887.448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2100.5
$ws.Range("I81").Value = 972
$ws.Range("K81").Value = 1944
$ws.Range("M81").Value = -883

$ws.Range("H84").Value = 2100.5
$ws.Range("I84").Value = 972
$ws.Range("K84").Value = 9720
$ws.Range("M84").Value = -4416

$ws.Range("H87").Value = 428627140
$ws.Range("I87").Value = 90000
$ws.Range("J87").Value = 500050000
$ws.Range("K87").Value = 90000
$ws.Range("L87").Value = 500050000
$ws.Range("M87").Value = -88752
$ws.Range("N87").Value = -500052496

$ws.Range("H90").Value = 428627140
$ws.Range("I90").Value = 90000
$ws.Range("J90").Value = 500050000
$ws.Range("K90").Value = 270000
$ws.Range("L90").Value = 1500150000
$ws.Range("M90").Value = -263760
$ws.Range("N90").Value = -1500162480

$ws.Range("H107").Value = 1815.7778
$ws.Range("J107").Value = 2700
$ws.Range("L107").Value = 8100
$ws.Range("N107").Value = -11940

$ws.Range("H122").Value = 5053.9
$ws.Range("I122").Value = 3106.5
$ws.Range("K122").Value = 9319.5
$ws.Range("M122").Value = -6869.5

$ws.Range("H132").Value = 4623.8237
$ws.Range("I132").Value = 4341.6665
$ws.Range("K132").Value = 13024.9995
$ws.Range("M132").Value = -10494.9995

$ws.Range("H136").Value = 2137.0435
$ws.Range("I136").Value = 1487.4667
$ws.Range("K136").Value = 4462.4001
$ws.Range("M136").Value = -1912.4001
